# "update factory for diagonal"
# Adds the "A* Diagonal" column (H) to the two "factory" result tables on
# Sheet1 (the header block at rows 29-35 and the header block at rows 37-50),
# mirroring the "A* Diagonal" column that already exists for the other
# blocks above (rows 1-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- header row 29 (factory map block) ---------------------------------
$ws.Range("H29").Value = "A* Diagonal"

# --- factory map data rows 38-42 (new column, previously absent) -------
$ws.Range("H38").Value = 3947
$ws.Range("H39").Value = 113.9
$ws.Range("H40").Value = 93
$ws.Range("H41").Value = 106
$ws.Range("H42").Value = 675

# --- factory terrain map data rows 45-49 (column existed but was empty) -
$ws.Range("H45").Value = 3946
$ws.Range("H46").Value = 117.6
$ws.Range("H47").Value = 93
$ws.Range("H48").Value = 110
$ws.Range("H49").Value = 540

# --- header row 50 (factory terrain map block) --------------------------
$ws.Range("H50").Value = "A* Diagonal"

# --- update the saved selection / scroll position -----------------------
$ws.Range("K41").Select()
